$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.187.48"
$ws.Range("D3").Value = "3.402.89"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.59"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.69"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("E7").Value = "  +7.41%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "3.403.57"
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.12"
$ws.Range("E10").Value = "  -2.97%  "
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.439"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").Value = "3.986.76"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000187"
$ws.Range("E15").Value = "  -3.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.63"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "64.208.22"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").Value = "3.451.22"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.31"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.84"
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.47"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.99"
$ws.Range("E22").Value = "  -2.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.546"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.78"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("E26").Value = "  -3.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.42"
$ws.Range("E27").Value = "  +6.43%  "
$ws.Range("E28").Value = "  -1.50%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.47"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.03"
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("E35").Value = "  +5.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.68"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.89"
$ws.Range("E38").Value = "  +5.28%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0760"
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").Value = "2.875.12"
$ws.Range("E40").Value = "  -4.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.35"
$ws.Range("E41").Value = "  -3.25%  "
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.84"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.29"
$ws.Range("E44").Value = "  +6.43%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0315"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.766"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "320.03"
$ws.Range("E47").Value = "  +6.05%  "
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.19"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.55"
$ws.Range("E51").Value = "  -0.98%  "
